# Updating attendance files 19th feb
# - Clear the "Day 8" header/lookup column (N4, N5, N6)
# - Flip individual P/A attendance marks for "Day 8" (column N) on several rows
# - Flip one additional mark in column L (row 23)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day-8 lookup formulas (N4/N5) and the day-number header (N6) are cleared entirely.
$ws.Range("N4").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("N6").ClearContents()

# Row 23 also has an L-column (Day 6) mark flip from P to A.
$ws.Range("L23").Value = "A"

# Column N (Day 8) attendance mark changes, row by row.
$nChanges = @{
    7  = "P"
    16 = "P"
    20 = "A"
    21 = "P"
    23 = "P"
    27 = "P"
    28 = "P"
    34 = "A"
    36 = "P"
    38 = "A"
    40 = "P"
    41 = "P"
    46 = "P"
    49 = "P"
    51 = "P"
    54 = "A"
    56 = "P"
    59 = "P"
    65 = "A"
    68 = "P"
    73 = "P"
}

foreach ($row in $nChanges.Keys) {
    $ws.Range("N$row").Value = $nChanges[$row]
}
